$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.602996254681648
$wsSummary.Range("C2").Value = 0.5625
$wsSummary.Range("D2").Value = 0.9269662921348315
$wsSummary.Range("E2").Value = 0.7001414427157001
$wsSummary.Range("F2").Value = 0.8206233421750663
$wsSummary.Range("G2").Value = 0.9044272663387211
$wsSummary.Range("H2").Value = 0.7682777146544348
$wsSummary.Range("I2").Value = 495
$wsSummary.Range("J2").Value = 385
$wsSummary.Range("K2").Value = 149
$wsSummary.Range("L2").Value = 39

# ---- Classification Report sheet ----
$wsReport = $wb.Worksheets.Item("Classification Report")

# row 2 - class "0"
$wsReport.Range("B2").Value = 0.7925531914893617
$wsReport.Range("C2").Value = 0.2790262172284644
$wsReport.Range("D2").Value = 0.4127423822714681

# row 3 - class "1"
$wsReport.Range("B3").Value = 0.5625
$wsReport.Range("C3").Value = 0.9269662921348315
$wsReport.Range("D3").Value = 0.7001414427157001

# row 4 - accuracy
$wsReport.Range("B4").Value = 0.602996254681648
$wsReport.Range("C4").Value = 0.602996254681648
$wsReport.Range("D4").Value = 0.602996254681648
$wsReport.Range("E4").Value = 0.602996254681648

# row 5 - macro avg
$wsReport.Range("B5").Value = 0.6775265957446808
$wsReport.Range("C5").Value = 0.602996254681648
$wsReport.Range("D5").Value = 0.5564419124935841

# row 6 - weighted avg
$wsReport.Range("B6").Value = 0.6775265957446808
$wsReport.Range("C6").Value = 0.602996254681648
$wsReport.Range("D6").Value = 0.5564419124935841

# ---- Confusion Matrix sheet ----
$wsConf = $wb.Worksheets.Item("Confusion Matrix")

# row 2 - Actual 0
$wsConf.Range("B2").Value = 149
$wsConf.Range("C2").Value = 385

# row 3 - Actual 1
$wsConf.Range("B3").Value = 39
$wsConf.Range("C3").Value = 495
